# Weinbar.xlsx edit: "added 1112 and freixenet"
#
# 1. Rename the existing "112 Blanc de Noirs" wine (row 100) to "1112 Blanc de Noirs".
# 2. Append a new wine row for "1112 Grauburgunder" (row 116).
# 3. Append a new wine row for "Freixenet Carta Rose Dry" (row 117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- 1. Fix the typo'd name in the existing row --
$ws.Range("A100").Value = "1112 Blanc de Noirs"

# -- 2. New row 116: 1112 Grauburgunder --
$ws.Range("A116").Value = "1112 Grauburgunder"
$ws.Range("B116").Value = "Grauburgunder"
$ws.Range("C116").Value = "white"
$ws.Range("D116").Value = "Germany"
$ws.Range("E116").Value = "Baden"
$ws.Range("F116").Value = "Markgraeflich Badisches Weinhaus"
$ws.Range("G116").Value = "suess, sueffig"
$ws.Range("H116").Value = "Rewe"
$ws.Range("I116").Value = 13
$ws.Range("J116").Value = 4.9
$ws.Range("K116").Value = 2019
$ws.Range("L116").Value = "yes"
$ws.Range("M116").Value = "yes"

# -- 3. New row 117: Freixenet Carta Rose Dry --
$ws.Range("A117").Value = "Freixenet Carta Rose Dry"
$ws.Range("B117").Value = "Garnacha, Trepat"
$ws.Range("C117").Value = "rose"
$ws.Range("D117").Value = "Spain"
$ws.Range("E117").Value = "Cava"
$ws.Range("F117").Value = "Henkell-Freixenet"
$ws.Range("G117").Value = "suess"
$ws.Range("H117").Value = "Rewe"
$ws.Range("I117").Value = 12
$ws.Range("J117").Value = 5.99
$ws.Range("K117").Value = 2019
$ws.Range("L117").Value = "yes"
$ws.Range("M117").Value = "yes"
